{"js": "// Update the ADHD.ASD / ASD summary-statistics column and their\n// dependent test statistics after excluding a participant who took\n// part a second time following a name change.\nconst replacements = [\n  [\"46.43 \u00b10.41 (32 to 67)\", \"46.64 \u00b10.43 (32 to 67)\"],\n  [\"32.54 \u00b10.33 (14 to 52)\", \"32.83 \u00b10.35 (14 to 52)\"],\n  [\"20.618*\", \"19.997*\"],\n  [\"30.22 \u00b10.36 (17 to 44)\", \"30.68 \u00b10.37 (17 to 44)\"],\n  [\"28.58 \u00b10.30 (18 to 45)\", \"29.04 \u00b10.30 (20 to 45)\"],\n  [\"-1.726\", \"-1.239\"],\n  [\"3.67 \u00b10.04 (2 to 5)\", \"3.65 \u00b10.04 (2 to 5)\"],\n  [\"0.226\", \"0.196\"],\n  [\"12 - 8 - 3\", \"12 - 7 - 3\"],\n  [\"12 - 12 - 0\", \"12 - 11 - 0\"],\n  [\"-4.406\", \"-4.021\"],\n  [\"112.93 \u00b10.49 (91 to 133)\", \"113.23 \u00b10.52 (91 to 133)\"],\n  [\"111.31 \u00b10.60 (78 to 144)\", \"111.98 \u00b10.63 (78 to 144)\"],\n  [\"-2.170\", \"-2.014\"],\n  [\"146.52 \u00b11.45 (85 to 201)\", \"145.86 \u00b11.55 (85 to 201)\"],\n  [\"152.92 \u00b11.70 (55 to 217)\", \"151.61 \u00b11.79 (55 to 217)\"],\n  [\"30.978*\", \"29.577*\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the ADHD.ASD / ASD summary-statistics column and their\n# dependent test statistics after excluding a participant who took\n# part a second time following a name change.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"46.43 \u00b10.41 (32 to 67)\", \"46.64 \u00b10.43 (32 to 67)\"),\n  @(\"32.54 \u00b10.33 (14 to 52)\", \"32.83 \u00b10.35 (14 to 52)\"),\n  @(\"20.618*\", \"19.997*\"),\n  @(\"30.22 \u00b10.36 (17 to 44)\", \"30.68 \u00b10.37 (17 to 44)\"),\n  @(\"28.58 \u00b10.30 (18 to 45)\", \"29.04 \u00b10.30 (20 to 45)\"),\n  @(\"-1.726\", \"-1.239\"),\n  @(\"3.67 \u00b10.04 (2 to 5)\", \"3.65 \u00b10.04 (2 to 5)\"),\n  @(\"0.226\", \"0.196\"),\n  @(\"12 - 8 - 3\", \"12 - 7 - 3\"),\n  @(\"12 - 12 - 0\", \"12 - 11 - 0\"),\n  @(\"-4.406\", \"-4.021\"),\n  @(\"112.93 \u00b10.49 (91 to 133)\", \"113.23 \u00b10.52 (91 to 133)\"),\n  @(\"111.31 \u00b10.60 (78 to 144)\", \"111.98 \u00b10.63 (78 to 144)\"),\n  @(\"-2.170\", \"-2.014\"),\n  @(\"146.52 \u00b11.45 (85 to 201)\", \"145.86 \u00b11.55 (85 to 201)\"),\n  @(\"152.92 \u00b11.70 (55 to 217)\", \"151.61 \u00b11.79 (55 to 217)\"),\n  @(\"30.978*\", \"29.577*\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Forward = $true\n  $find.Wrap = 1        # wdFindContinue\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.MatchSoundsLike = $false\n  $find.MatchAllWordForms = $false\n\n  $found = $find.Execute(\n    $find.Text,\n    $find.MatchCase,\n    $find.MatchWholeWord,\n    $find.MatchWildcards,\n    $find.MatchSoundsLike,\n    $find.MatchAllWordForms,\n    $find.Forward,\n    $find.Wrap,\n    $find.Format,\n    $find.Replacement.Text,\n    2                    # wdReplaceAll\n  )\n\n  if (-not $found) {\n    throw \"Text not found: $old\"\n  }\n}\n"}
